# This script applies the "Updated symbol list" commit: a batch of price /
# volume-label refreshes plus a 3-row re-ordering (KickToken, BKEXToken,
# CEJI) in the cryptos sheet.
#
# All of the touched columns (B, C, D, E) store their data as TEXT, even
# though many of the values look numeric (e.g. "249.34"). Setting
# Range.Value directly with a numeric-looking string causes Excel to coerce
# it into a real number, which would change the cell's underlying type and
# break the diff. To avoid that we:
#   1. Prefix the value with a single quote so Excel treats it as text
#      ("quote-prefixed" literal), and
#   2. Reset the range's Style back to "Normal" afterwards, which clears the
#      quote-prefix formatting flag without touching the text value, so the
#      resulting cell looks just like an ordinary text cell (no extra style
#      applied), matching the original workbook's formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $range = $ws.Range($Address)
    $range.Value = "'" + $Text
    $range.Style = "Normal"
}

# --- Simple price / label updates -----------------------------------------
Set-TextValue "D2"  "249.34"
Set-TextValue "D3"  "22.93"
Set-TextValue "D4"  "5.433"
Set-TextValue "D5"  "0.05623"
Set-TextValue "D7"  "6.369"
Set-TextValue "D8"  "0.8114"
Set-TextValue "D9"  "0.8977"
Set-TextValue "D10" "0.1429"
Set-TextValue "D11" "0.07514"
Set-TextValue "D12" "0.03096"
Set-TextValue "D13" "0.03091"
Set-TextValue "D14" "0.09324"
Set-TextValue "D15" "3.569"
Set-TextValue "D16" "0.001594"
Set-TextValue "E18" "17OneONEWorstin24h"
Set-TextValue "D19" "0.006416"
Set-TextValue "D20" "0.004991"
Set-TextValue "D21" "0.001034"
Set-TextValue "D23" "3.703"
Set-TextValue "D24" "2.180"
Set-TextValue "D25" "0.3302"
Set-TextValue "D40" "0.04053"

# --- Rows 41-43 get re-ordered (Kick / BKEX / CEJI) and re-valued ---------
Set-TextValue "B41" "KickToken"
Set-TextValue "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006790"
Set-TextValue "E41" "40KickTokenKICK"

Set-TextValue "B42" "BKEXToken"
Set-TextValue "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1067"
Set-TextValue "E42" "41BKEXTokenBKK"

Set-TextValue "B43" "CEJI"
Set-TextValue "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002736"
Set-TextValue "E43" "42CEJICEJI"

# --- More simple price updates ---------------------------------------------
Set-TextValue "D44" "0.007474"
Set-TextValue "D45" "0.00005578"
Set-TextValue "D47" "0.5005"
Set-TextValue "D48" "0.2400"
